# Adding Custom Interactable table
# The shot-tracking table on Sheet1 gets two cell updates:
#   - G6 (ARTIST for RL_03_SH_0260): "Vishal Mahale" -> "Rohit Chavan"
#   - J8 (NOTES for RL_02_SH_0190): "Roto the fg character" -> "Roto all characters"
# and the active selection moves from E2 to G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Rohit Chavan"
$ws.Range("J8").Value = "Roto all characters"

$ws.Range("G11").Select()
